$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "hehe"
$ws.Range("B2").Value = "Item 1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Desc 2"
$ws.Range("F2").Value = "Completed"

# Row 3
$ws.Range("A3").Value = "hehe"
$ws.Range("B3").Value = "Item 111111111111111111111111111111111111111111"
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = "Desc 3"
$ws.Range("F3").Value = "Completed"
